$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Build the Excel Table ("Table1") covering A1:U55, with renamed headers:
#       "<name>_old" -> "<name>_FV2404"
#       "<name>_new" -> "<name>_FV2410"
#    (column K / "diff" keeps its name).
#
#    The table is first created on a throw-away blank range so that the
#    engine doesn't have to fold the *existing* (bold/shaded) header-row
#    formatting into a new dxf - then it is resized onto the real header
#    row and every column is renamed in place (which also rewrites the
#    actual header cell text, reusing the pre-existing cell style).
# ---------------------------------------------------------------------------
$scratch = $ws.Range("W1:W2")
$scratch.Cells.Item(1, 1).Value = "scratch_header"

$tbl = $ws.ListObjects.Add(1, $scratch, $null, 1)
$tbl.Name = "Table1"
$tbl.Resize($ws.Range("A1:U55"))

$columnNames = @(
    "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404", "Segment ID_FV2404",
    "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404", "Bedingungsausdruck_FV2404", "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410", "Segment ID_FV2410",
    "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410", "Bedingungsausdruck_FV2410", "Bedingung_FV2410"
)

for ($i = 1; $i -le $columnNames.Count; $i++) {
    $tbl.ListColumns.Item($i).Range.Cells.Item(1, 1).Value = $columnNames[$i - 1]
}

# Drop the scratch cells used to bootstrap the table so they don't linger
# outside the real A1:U55 range.
$ws.Range("W1:W2").Clear()

# ---------------------------------------------------------------------------
# 2) Freeze the header row (pane split after row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
